# Feature/du/update template for active learning
# - Remove the "Currency" / "USD" column (column S) from the "Simple Fields"
#   and "Simple Fields - Formatted" sheets, shifting "Items"/"table" left.
# - Change the Payment Terms value from "due 30 days" to "30 days" on both
#   of those sheets.
# The "Items" and "Items - Formatted" sheets keep identical content - only
# the shared-string bookkeeping shifts as a side effect of removing the
# "Currency"/"USD" strings, so nothing needs to change there explicitly.

$wb = $excel.ActiveWorkbook

$simpleSheets = @("Simple Fields", "Simple Fields - Formatted")

foreach ($sheetName in $simpleSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Update the Payment Terms text in row 2 (column M) before the shift.
    $ws.Range("M2").Value = "30 days"

    # Delete the whole "Currency" column (S); this shifts column T ("Items")
    # left into column S, matching the target layout (A1:S2).
    $ws.Range("S1").EntireColumn.Delete()
}
